$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1113
$ws1.Range("F6").Value = 617
$ws1.Range("F9").Value = 141
$ws1.Range("F10").Value = 141
$ws1.Range("F11").Value = 1431
$ws1.Range("F12").Value = 3060
$ws1.Range("F13").Value = 581
$ws1.Range("F14").Value = 1732
$ws1.Range("F15").Value = 1788
$ws1.Range("F16").Value = 833
$ws1.Range("F17").Value = 265
$ws1.Range("F22").Value = 392
$ws1.Range("F25").Value = 4696
$ws1.Range("F26").Value = 743
$ws1.Range("F28").Value = 1620

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 25
$ws2.Range("F5").Value = 22
$ws2.Range("F6").Value = 45

# Sheet "本地生活" (local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 33

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 33
$ws4.Range("F5").Value = 25
$ws4.Range("F8").Value = 22
$ws4.Range("F9").Value = 45
$ws4.Range("F12").Value = 1113
$ws4.Range("F14").Value = 617
$ws4.Range("F17").Value = 141
$ws4.Range("F18").Value = 141
$ws4.Range("F20").Value = 1431
$ws4.Range("F21").Value = 3060
$ws4.Range("F22").Value = 581
$ws4.Range("F23").Value = 1732
$ws4.Range("F24").Value = 1788
$ws4.Range("F25").Value = 833
$ws4.Range("F26").Value = 265
$ws4.Range("F33").Value = 392
$ws4.Range("F36").Value = 4696
$ws4.Range("F37").Value = 743
$ws4.Range("F39").Value = 1620
